$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force target cells to text format so values are not auto-converted
# to dates/numbers (the source data must remain literal text).
$ws.Range("A1:D1").NumberFormat = "@"
$ws.Range("A2:E10").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "Date collected"
$ws.Range("B1").Value = "Plot__"
$ws.Range("C1").Value = "Species"
$ws.Range("D1").Value = "Sex Weight"

# Data rows
$data = @(
    @("1/9/78",   "1",  "DM", "M",  "40"),
    @("1/9/78",   "1",  "DM", "EF", "36"),
    @("1/9/78",   "41", "DS", "EF", "135"),
    @("1/20/78",  "1",  "DM", "F",  "39"),
    @("1/20/78",  "2",  "DM", "M",  "43"),
    @("1/20/78",  "2",  "DS", "EF", "144"),
    @("3/13/78",  "2",  "DM", "EF", "51"),
    @("3/13/78",  "2",  "DM", "EF", "44"),
    @("3/13/78",  "2",  "DS", "EF", "146")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
